$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I, J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header-row formatting (bold, bordered, centered)
# by copying the style from the adjacent header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for column I (I0) rows 2-30
$iValues = @(5,8,6,9,8,8,8,4,9,8,9,7,8,3,7,6,6,8,8,6,9,6,7,6,5,6,4,7,6)
# Data values for column J (IF) rows 2-30
$jValues = @(7,8,8,9,8,8,8,6,9,9,9,8,9,5,7,7,6,8,8,6,9,6,7,7,5,7,5,7,7)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
